$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2700
$ws.Range("I111").Value = 1812.5
$ws.Range("J111").Value = 6250
$ws.Range("K111").Value = 5437.5
$ws.Range("L111").Value = 18750
$ws.Range("M111").Value = -2370.5
$ws.Range("N111").Value = -24884
$ws.Range("H132").Value = 2003.7759
$ws.Range("I132").Value = 1210.14
$ws.Range("J132").Value = 6964
$ws.Range("K132").Value = 3630.42
$ws.Range("L132").Value = 20892
$ws.Range("M132").Value = -1100.42
$ws.Range("N132").Value = -25952
$ws.Range("H138").Value = 2364.6404
$ws.Range("I138").Value = 2247.074
$ws.Range("J138").Value = 2415.8386
$ws.Range("K138").Value = 6741.222
$ws.Range("L138").Value = 7247.5158
$ws.Range("M138").Value = -1601.222
$ws.Range("N138").Value = -17527.5158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21324.988
$ws.Range("I32").Value = 5878.3975
$ws.Range("J32").Value = 141808.4
$ws.Range("K32").Value = 5878.3975
$ws.Range("L32").Value = 141808.4
$ws.Range("M32").Value = -5591.3975
$ws.Range("N32").Value = -142382.4
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").ClearContents()
$ws.Range("N44").Value = 0
$ws.Range("H55").Value = 37000
$ws.Range("J55").Value = 37000
$ws.Range("L55").Value = 37000
$ws.Range("N55").Value = -37630
$ws.Range("H74").Value = 1790.2683
$ws.Range("I74").Value = 717.17145
$ws.Range("K74").Value = 717.17145
$ws.Range("M74").Value = 156.82855
$ws.Range("H77").Value = 1790.2683
$ws.Range("I77").Value = 717.17145
$ws.Range("K77").Value = 3585.85725
$ws.Range("M77").Value = 782.14275
$ws.Range("H80").Value = 38000
$ws.Range("J80").Value = 38000
$ws.Range("L80").Value = 38000
$ws.Range("N80").Value = -39996
$ws.Range("H83").Value = 38000
$ws.Range("J83").Value = 38000
$ws.Range("L83").Value = 114000
$ws.Range("N83").Value = -123984

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 46685.668
$ws.Range("J9").Value = 46685.668
$ws.Range("L9").Value = 46685.668
$ws.Range("N9").Value = -47021.668
$ws.Range("H35").Value = 30400
$ws.Range("J35").Value = 30400
$ws.Range("L35").Value = 30400
$ws.Range("N35").Value = -31020
$ws.Range("H82").Value = 21923.188
$ws.Range("I82").Value = 10793
$ws.Range("J82").Value = 30580
$ws.Range("K82").Value = 10793
$ws.Range("L82").Value = 30580
$ws.Range("M82").Value = -10410
$ws.Range("N82").Value = -31346
$ws.Range("H85").Value = 21923.188
$ws.Range("I85").Value = 10793
$ws.Range("J85").Value = 30580
$ws.Range("K85").Value = 10793
$ws.Range("L85").Value = 30580
$ws.Range("M85").Value = -9467
$ws.Range("N85").Value = -33232
$ws.Range("H133").Value = 32000
$ws.Range("J133").Value = 32000
$ws.Range("L133").Value = 32000
$ws.Range("N133").Value = -42120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3155.7666
$ws.Range("I31").Value = 1416.32
$ws.Range("J31").Value = 4398.2285
$ws.Range("K31").Value = 1416.32
$ws.Range("L31").Value = 4398.2285
$ws.Range("M31").Value = -1121.32
$ws.Range("N31").Value = -4988.2285
$ws.Range("H34").Value = 3155.7666
$ws.Range("I34").Value = 1416.32
$ws.Range("J34").Value = 4398.2285
$ws.Range("K34").Value = 1416.32
$ws.Range("L34").Value = 4398.2285
$ws.Range("M34").Value = -1214.32
$ws.Range("N34").Value = -4802.2285
$ws.Range("H41").Value = 12741.25
$ws.Range("J41").Value = 21532.5
$ws.Range("L41").Value = 21532.5
$ws.Range("N41").Value = -22388.5
$ws.Range("H50").Value = 9223.6
$ws.Range("J50").Value = 9223.6
$ws.Range("L50").Value = 9223.6
$ws.Range("N50").Value = -10473.6
$ws.Range("H60").Value = 23701.533
$ws.Range("J60").Value = 23701.533
$ws.Range("L60").Value = 23701.533
$ws.Range("N60").Value = -24723.533
$ws.Range("H68").Value = 17501.666
$ws.Range("J68").Value = 17501.666
$ws.Range("L68").Value = 17501.666
$ws.Range("N68").Value = -18999.666
$ws.Range("H71").Value = 17501.666
$ws.Range("J71").Value = 17501.666
$ws.Range("L71").Value = 52504.99800000001
$ws.Range("N71").Value = -59992.99800000001
$ws.Range("H94").Value = 4298.8823
$ws.Range("I94").Value = 925
$ws.Range("J94").Value = 4748.7334
$ws.Range("K94").Value = 925
$ws.Range("L94").Value = 4748.7334
$ws.Range("M94").Value = -474
$ws.Range("N94").Value = -5650.7334
$ws.Range("H109").Value = 11950
$ws.Range("J109").Value = 11950
$ws.Range("L109").Value = 11950
$ws.Range("N109").Value = -14030

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 18371.834
$ws.Range("J57").Value = 25030.25
$ws.Range("L57").Value = 25030.25
$ws.Range("N57").Value = -26670.25
$ws.Range("H102").Value = 1411.0714
$ws.Range("I102").Value = 1034.6086
$ws.Range("J102").Value = 3142.8
$ws.Range("K102").Value = 1034.6086
$ws.Range("L102").Value = 3142.8
$ws.Range("M102").Value = 587.3914
$ws.Range("N102").Value = -6386.8
$ws.Range("H123").Value = 38973.25
$ws.Range("J123").Value = 38973.25
$ws.Range("L123").Value = 38973.25
$ws.Range("N123").Value = -43873.25
$ws.Range("H132").Value = 2828.4583
$ws.Range("I132").Value = 2605.6667
$ws.Range("J132").Value = 3199.7778
$ws.Range("K132").Value = 7817.000100000001
$ws.Range("L132").Value = 9599.3334
$ws.Range("M132").Value = -5287.000100000001
$ws.Range("N132").Value = -14659.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1879.1666
$ws.Range("I46").Value = 1647.3684
$ws.Range("J46").Value = 2760
$ws.Range("K46").Value = 1647.3684
$ws.Range("L46").Value = 2760
$ws.Range("M46").Value = -1459.3684
$ws.Range("N46").Value = -3136

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 2991.158
$ws.Range("I107").Value = 1543.4
$ws.Range("J107").Value = 4599.778
$ws.Range("K107").Value = 4630.200000000001
$ws.Range("L107").Value = 13799.334
$ws.Range("M107").Value = -2710.200000000001
$ws.Range("N107").Value = -17639.334
$ws.Range("H109").Value = 27800
$ws.Range("J109").Value = 27800
$ws.Range("L109").Value = 27800
$ws.Range("N109").Value = -30574
$ws.Range("H122").Value = 9215.296
$ws.Range("I122").Value = 11735.526
$ws.Range("K122").Value = 35206.578
$ws.Range("M122").Value = -32756.578
